$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 8: "All" contrast summing the first five individual contrasts
$ws.Range("A8").Value = "All"
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").HorizontalAlignment = -4131  # xlLeft

$values = @(1, 1, 1, 1, 1, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 2 + $i  # Column B = 2
    $cell = $ws.Cells.Item(8, $col)
    $cell.Value = $values[$i]
    $cell.HorizontalAlignment = -4131  # xlLeft
}

# Update the active selection to match the new last-used cell
$ws.Range("L8").Select()
